# Trade #108 closed at 2026-02-16 21:42:39 - leadlag DOWN +0.000%
#
# This script:
#  1) Updates the Summary sheet's aggregate stats (rows 2 & 3)
#  2) Closes two previously-OPEN leadlag trades (rows 63 & 64) and
#     appends a brand-new OPEN leadlag trade (row 84) to the "leadlag" sheet
#  3) Appends the two newly-closed trades to the "All Trades" sheet (rows 84 & 85)
#  4) Updates the Comparison sheet's leadlag summary row (row 2)

$wb = $excel.ActiveWorkbook

# Helper: write a literal text value into a cell, forcing Excel to treat it
# as plain text rather than auto-converting it into a number/date/percentage.
function Set-CellText($sheet, $addr, [string]$text) {
    $sheet.Range($addr).Value = "'" + $text
}

# Helper: write a numeric value into a cell.
function Set-CellNumber($sheet, $addr, $num) {
    $sheet.Range($addr).Value = $num
}

# ---------------------------------------------------------------------
# 1) Summary sheet
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")

Set-CellNumber $wsSummary "C2" 84
Set-CellText   $wsSummary "D2" "70.2%"
Set-CellText   $wsSummary "E2" "+24.7469%"
Set-CellText   $wsSummary "F2" "+0.2946%"

Set-CellNumber $wsSummary "C3" 82
Set-CellText   $wsSummary "D3" "48.8%"
Set-CellText   $wsSummary "E3" "+13.4256%"
Set-CellText   $wsSummary "F3" "+0.1637%"

# ---------------------------------------------------------------------
# 2) leadlag sheet
# ---------------------------------------------------------------------
$wsLeadlag = $wb.Worksheets.Item("leadlag")

# Row 63: trade gets closed
Set-CellNumber $wsLeadlag "G63" 68719.82691600001
Set-CellText   $wsLeadlag "H63" "CLOSED"
Set-CellNumber $wsLeadlag "I63" 0.1284
Set-CellNumber $wsLeadlag "J63" 1.28
Set-CellText   $wsLeadlag "M63" "time_exit_5min"
Set-CellNumber $wsLeadlag "N63" 5

# Row 64: trade gets closed
Set-CellNumber $wsLeadlag "G64" 68492.04945400001
Set-CellText   $wsLeadlag "H64" "CLOSED"
Set-CellNumber $wsLeadlag "I64" 0.0992
Set-CellNumber $wsLeadlag "J64" 0.99
Set-CellText   $wsLeadlag "M64" "time_exit_5min"
Set-CellNumber $wsLeadlag "N64" 5

# Row 84: brand-new OPEN trade (#108) appended
Set-CellNumber $wsLeadlag "A84" 108
Set-CellText   $wsLeadlag "B84" "2026-02-16"
Set-CellText   $wsLeadlag "C84" "21:42:39"
Set-CellText   $wsLeadlag "D84" "leadlag"
Set-CellText   $wsLeadlag "E84" "DOWN"
Set-CellNumber $wsLeadlag "F84" 68401.755
# G84 stays blank (exit price not yet known - trade is OPEN)
Set-CellText   $wsLeadlag "H84" "OPEN"
Set-CellNumber $wsLeadlag "I84" 0
Set-CellNumber $wsLeadlag "J84" 0
Set-CellNumber $wsLeadlag "K84" 0.75
Set-CellText   $wsLeadlag "L84" "Binance leading with -0.103% move"
# M84 stays blank (exit reason not yet known - trade is OPEN)
Set-CellNumber $wsLeadlag "N84" 0

# ---------------------------------------------------------------------
# 3) All Trades sheet - append the two newly closed trades
# ---------------------------------------------------------------------
$wsAllTrades = $wb.Worksheets.Item("All Trades")

# Row 84 (mirrors leadlag row 63, now closed)
Set-CellNumber $wsAllTrades "A84" 83
Set-CellText   $wsAllTrades "B84" "2026-02-16"
Set-CellText   $wsAllTrades "C84" "21:37:19"
Set-CellText   $wsAllTrades "D84" "leadlag"
Set-CellText   $wsAllTrades "E84" "UP"
Set-CellNumber $wsAllTrades "F84" 68631.675
Set-CellNumber $wsAllTrades "G84" 68719.82691600001
Set-CellText   $wsAllTrades "H84" "CLOSED"
Set-CellNumber $wsAllTrades "I84" 0.1284
Set-CellNumber $wsAllTrades "J84" 1.28
Set-CellNumber $wsAllTrades "K84" 0.75
Set-CellText   $wsAllTrades "L84" "Coinbase leading with 0.114% move"
Set-CellText   $wsAllTrades "M84" "time_exit_5min"
Set-CellNumber $wsAllTrades "N84" 5

# Row 85 (mirrors leadlag row 64, now closed)
Set-CellNumber $wsAllTrades "A85" 84
Set-CellText   $wsAllTrades "B85" "2026-02-16"
Set-CellText   $wsAllTrades "C85" "21:37:36"
Set-CellText   $wsAllTrades "D85" "leadlag"
Set-CellText   $wsAllTrades "E85" "DOWN"
Set-CellNumber $wsAllTrades "F85" 68560.03999999999
Set-CellNumber $wsAllTrades "G85" 68492.04945400001
Set-CellText   $wsAllTrades "H85" "CLOSED"
Set-CellNumber $wsAllTrades "I85" 0.0992
Set-CellNumber $wsAllTrades "J85" 0.99
Set-CellNumber $wsAllTrades "K85" 0.75
Set-CellText   $wsAllTrades "L85" "Binance leading with -0.111% move"
Set-CellText   $wsAllTrades "M85" "time_exit_5min"
Set-CellNumber $wsAllTrades "N85" 5

# ---------------------------------------------------------------------
# 4) Comparison sheet
# ---------------------------------------------------------------------
$wsComparison = $wb.Worksheets.Item("Comparison")

Set-CellNumber $wsComparison "B2" 82
Set-CellText   $wsComparison "C2" "48.8%"
Set-CellText   $wsComparison "D2" "2.84"
Set-CellText   $wsComparison "E2" "+0.5177%"
Set-CellText   $wsComparison "G2" "1.63"

Write-Host "Edit complete."
